$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A4").Value = "batch_003"
$ws.Range("B4").Value = "y"
$ws.Range("C4").Value = "批量操作语句3执行"
$ws.Range("D4").Value = "batchsql"
$ws.Range("F4").Value = "batch03"
$ws.Range("H4").Value = "batch_sql_03"
$ws.Range("I4").Value = "select * from `$batch03"
$ws.Range("J4").Value = "src/test/resources/io.dingodb.test/testdata/cases/batchsql/expectedresult/batch_003.csv"
$ws.Range("M4").Value = "csv_containsAll"

$ws.Range("M8").Select()
$ws.Application.ActiveWindow.ScrollColumn = 2
